$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Veda")

# Update the formula in Q10 to include the R26/R27 adjustment
$ws.Range("Q10").Formula = "=SUMIFS(iea_data!I3:I9999,iea_data!`$B`$3:`$B`$9999,Veda!`$Q`$9)+R26-R27"

# Recalculate the workbook so dependent formulas/chart caches refresh
$excel.CalculateFullRebuild()

# Update the active selection to Q10 as in the edited workbook
$ws.Activate()
$ws.Range("Q10").Select()
